# Apply weekly fruit/vegetable price update: the D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) values are redistributed across rows 2-11 according to
# the mapping below (new row -> source row, using the ORIGINAL data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns we need to move as a group, per row
$cols = @("D", "J", "K", "L", "M", "P")

# new row number -> row number that currently holds the data that should end up there
$mapping = @{
    2  = 4
    3  = 5
    4  = 7
    5  = 6
    6  = 11
    7  = 10
    8  = 3
    9  = 2
    10 = 9
    11 = 8
}

# Snapshot the current (original) values for each row/column before any writes,
# since several rows both give and receive values (e.g. 2<->4, 9<->2, etc.)
$original = @{}
foreach ($row in 2..11) {
    $original[$row] = @{}
    foreach ($col in $cols) {
        $original[$row][$col] = $ws.Range("$col$row").Value2()
    }
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $original[$srcRow][$col]
    }
}
